$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) "Purpose" paragraph: split the single run into two runs - trim the
#    trailing space from the first run and append a new run with the new
#    sentence about the challenges section.
# ---------------------------------------------------------------------------
$purposePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Its purpose is to give chances*") {
        $purposePara = $cand
        break
    }
}

$purposeRange = $d.Range($purposePara.Range.Start, $purposePara.Range.End)
$purposeXml = "<w:p $wns>" +
    "<w:pPr><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t>Its purpose is to give chances to everyone willing to showcase their skills, creating content, they could also get feedbacks and motivation on their work through comments, and there is a section where viewers can donate to your account just to support your work and art</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t>, there is a section of challenges where you can do completions and win prices with your creativity and skills</w:t></w:r>" +
    "</w:p>"
$purposeRange.InsertXML($purposeXml)

# ---------------------------------------------------------------------------
# 2) Insert a new bullet paragraph ("Locating (...)") right before the
#    "Problems" paragraph.
# ---------------------------------------------------------------------------
$problemsPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "Problems`r") {
        $problemsPara = $cand
        break
    }
}

$insertPoint = $d.Range($problemsPara.Range.Start, $problemsPara.Range.Start)
$locatingXml = "<w:p $wns>" +
    "<w:pPr><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:i/><w:iCs/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t>Locating (</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t>you can locate your nearest photographer in your area, you DM and book your session)</w:t></w:r>" +
    "</w:p>"
$insertPoint.InsertXML($locatingXml)

# ---------------------------------------------------------------------------
# 3) Move the <w:lastRenderedPageBreak/> marker from the lone-space paragraph
#    (just after "Lack of sleep...") onto the "Lacking of mobile data
#    bundles..." paragraph.
# ---------------------------------------------------------------------------
$lackingPara = $null
$spacePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Lacking of mobile data bundles*") {
        $lackingPara = $cand
    }
    if ($cand.Range.Text -eq " `r") {
        $spacePara = $cand
    }
}

$lackingRange = $d.Range($lackingPara.Range.Start, $lackingPara.Range.End)
$lackingXml = "<w:p $wns>" +
    "<w:pPr><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:lastRenderedPageBreak/><w:t>Lacking of mobile data bundles, had to struggle to get data to search for other information on the internet.</w:t></w:r>" +
    "</w:p>"
$lackingRange.InsertXML($lackingXml)

$spaceRange = $d.Range($spacePara.Range.Start, $spacePara.Range.End)
$spaceXml = "<w:p $wns>" +
    "<w:pPr><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "</w:p>"
$spaceRange.InsertXML($spaceXml)

Write-Output "done"
